# Update handback status timestamps ("Generate Report for Handback")

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 326a0a1a... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-19 20:50:21"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 326a0a1a... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-19 20:50:17"
$wsZhCn.Range("K3").Value = "2016-08-19 20:50:33"

# de-de sheet: Correspond Handback DateTime for 326a0a1a... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-19 20:50:40"
